# Append the 2021/10/11 DGS report row to the risk-matrix time series.
#
# New row 92: date "2021/10/11", incidencia_portugal=82.9,
# incidencia_continente=82.7, r_portugal=0.95, r_continente=0.95

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 92

# Column A holds the date as literal text (e.g. "2021/10/08") even though the
# cell carries a yyyy/mm/dd number format - that is how every other row in
# this sheet stores its date. Assigning a date-shaped string straight into a
# date-formatted cell via .Value/.Formula gets auto-parsed into a real date
# serial, which would not match the existing rows. To avoid that, build the
# text on an out-of-the-way scratch cell via a text formula (so it is never
# auto-converted), copy it, and paste just the resulting value into the
# target cell - this preserves both the "stored as text" shared-string type
# and the existing date style/format of the column.
$scratchRow = 1000
$ws.Cells.Item($scratchRow, 1).Formula = '="2021/10/11"'
$ws.Cells.Item($scratchRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163) # xlPasteValues
$ws.Cells.Item($scratchRow, 1).Formula = ""
$ws.Rows.Item($scratchRow).Delete()

# Numeric columns
$ws.Cells.Item($newRow, 2).Value = 82.9
$ws.Cells.Item($newRow, 3).Value = 82.7
$ws.Cells.Item($newRow, 4).Value = 0.95
$ws.Cells.Item($newRow, 5).Value = 0.95

# Move the selection to the next empty row, mirroring where Excel would leave
# the cursor after the new row was entered.
$ws.Range("A93").Select()
